# TT-moduly.xlsx: name modules 36 and 37 (from the "kroužek"/circle batch),
# record the production-start date for module 37, and update the active
# selection to reflect where editing left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SeznamModulu")

# Module 36 (row 45): "Kroužek 13/14" -> "Bunkr (kroužek 13/14)"
$ws.Range("B45").Value = "Bunkr (kroužek 13/14)"

# Module 37 (row 46): "Kroužek 13/14 II" -> "Tábor (kroužek 16/17)"
$ws.Range("B46").Value = "Tábor (kroužek 16/17)"

# Module 37 now has a recorded "Začátek výroby" (production start) date.
$ws.Range("C46").Value = "9/1/2016"

# Leave the selection where the edit was last made.
$ws.Range("C47").Select()
